# Listas-comparacao.xlsx update
# - Update the "Numero de elementos" / "Numero de nos consultados" figures for
#   the PLANO 1 comparison table (rows 9-10, columns E and G).
# - Widen columns E and G so the larger numbers keep fitting ("best fit").
# - Add a new styled (underlined) marker cell at H14, which extends the used
#   range and moves the active selection there, mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data values ---------------------------------------------
$ws.Range("E9").Value2  = 50005000
$ws.Range("G9").Value2  = 49985001
$ws.Range("E10").Value2 = 4501500
$ws.Range("G10").Value2 = 8994002

# --- Widen columns E and G to fit the new (larger) numbers --------------
$ws.Columns("E").ColumnWidth = 8
$ws.Columns("G").ColumnWidth = 8

# --- Add the new underlined marker cell at H14 ---------------------------
$ws.Range("H14").Font.Underline = 2

# --- Move the active selection to the newly added cell -------------------
$ws.Range("H14").Select() | Out-Null
